# Ajustes de actas y envio ftp
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Datos del usuario: nuevo responsable y correo
$ws.Range("C5").Value = "ALEXANDER CHRISTIAN FLORES CASTILLO"
$ws.Range("C6").Value = "aflorescast@pj.gob.pe"
# Pie de firma (mismo nombre del usuario, repetido bajo la linea de firma)
$ws.Range("F19").Value = "ALEXANDER CHRISTIAN FLORES CASTILLO"

# Direccion (se agrega el punto luego de "AV")
$ws.Range("C9").Value = "AV. SIGLO XX S/N"

# DNI (texto, conserva el formato original de la celda mediante comilla inicial)
$ws.Range("G5").Formula = "'76639137"

# Tipo de acta
$ws.Range("I4").Value = "DEVOLUCIÓN"

# Folio / secuencia de acta
$ws.Range("I7").Value = "1-2024"

# Fila del bien patrimonial (orden, codigo, denominacion, marca, modelo, serie)
# La celda A14 tiene formato de texto ("@") pero el valor se guarda como numero,
# igual que en el archivo original; forzamos el tipo numerico y luego restauramos
# el formato de la celda.
$ws.Range("A14").NumberFormat = "General"
$ws.Range("A14").Value = 1
$ws.Range("A14").NumberFormat = "@"
$ws.Range("B14").Formula = "'740800010005"
$ws.Range("C14").Value = "IMPRESORA DE CHEQUES"
$ws.Range("D14").Value = "DELL"
$ws.Range("E14").Value = "OptiPlex 3060"
$ws.Range("G14").Value = "HTHNBNN"
